# TC09_CDS_Filter_PHSAccession-phs002371.xlsx correction
# The FilesTab query (cell B4 on Sheet1) hard-coded
# experimental_strategies: ["RNA-Seq"] as a filter input. Per the commit
# ("obj correction input file correction") this should be reset back to an
# empty filter (no experimental-strategy restriction), matching the other
# three query cells on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(4, 2)   # B4 - FilesTab "query" column
$text = $cell.Value()

$updated = $text.Replace('experimental_strategies: ["RNA-Seq"]', 'experimental_strategies: []')
$cell.Value = $updated

# Keep the original row heights (they were auto-sized for the long query
# text; restore them explicitly since updating the text otherwise causes
# the rows to auto-fit to the new, shorter content).
$ws.Rows.Item(2).RowHeight = 409.5
$ws.Rows.Item(3).RowHeight = 409.5
$ws.Rows.Item(4).RowHeight = 409.5

# Reflect that B4 is the cell that was edited/selected.
$cell.Select()
